# Update formatting of sample cookies
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create / reuse a custom number format for datetime columns
$dateFormat = "YYYY-MM-DD HH:MM:SS"

# --- Row 2 ---
$ws.Range("A2").Value = 45656.01955134795
$ws.Range("A2").NumberFormat = $dateFormat

$ws.Range("H2").Value = 45657.80829861111
$ws.Range("H2").NumberFormat = $dateFormat

$ws.Range("I2").Value = $true
$ws.Range("J2").Value = $false
$ws.Range("K2").ClearContents()
$ws.Range("L2").Value = $true
$ws.Range("M2").Value = $true
$ws.Range("N2").Value = "Medium"
$ws.Range("O2").Value = "Lax"

$ws.Range("R2").Value = 45656.01955134718
$ws.Range("R2").NumberFormat = $dateFormat

$ws.Range("T2").Value = $true

# --- Row 3 ---
$ws.Range("A3").Value = 45656.01955137827
$ws.Range("A3").NumberFormat = $dateFormat

$ws.Range("H3").Value = 45657.80902777778
$ws.Range("H3").NumberFormat = $dateFormat

$ws.Range("I3").Value = $true
$ws.Range("J3").Value = $true
$ws.Range("K3").ClearContents()
$ws.Range("L3").Value = $true
$ws.Range("M3").Value = $true
$ws.Range("N3").Value = "Medium"
$ws.Range("O3").Value = "Strict"

$ws.Range("R3").Value = 45656.01955137773
$ws.Range("R3").NumberFormat = $dateFormat

$ws.Range("T3").Value = $true

# --- Row 4 ---
$ws.Range("A4").Value = 45656.01955140649
$ws.Range("A4").NumberFormat = $dateFormat

$ws.Range("H4").ClearContents()

$ws.Range("I4").Value = $true
$ws.Range("J4").Value = $true
$ws.Range("K4").ClearContents()
$ws.Range("L4").Value = $false
$ws.Range("M4").Value = $false
$ws.Range("N4").Value = "Medium"
$ws.Range("O4").Value = "None"

$ws.Range("R4").Value = 45656.01955140598
$ws.Range("R4").NumberFormat = $dateFormat

$ws.Range("T4").Value = $true
